$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.613956809043884
$ws.Range("B1").Value = 1.751061201095581
$ws.Range("C1").Value = 2.010695934295654
$ws.Range("D1").Value = 2.501955270767212
$ws.Range("E1").Value = 1.710631132125854
